function Set-TextValue($ws, $addr, $val) {
    # Force the cell to store $val as literal text (matches the source
    # data which is always written as inline/shared strings, never as
    # numbers) - plain .Value assignment would let Excel auto-coerce
    # numeric-looking strings (e.g. "0.7723", "244.52") into real numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    # Reset the style back to the default ("Normal") so only the cells
    # stored type changes - the visible/persisted style stays untouched,
    # matching the workbook where these data cells carry no explicit style.
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws "D2" "29.898.03"
Set-TextValue $ws "E2" "  +0.13%  "

# Row 3
Set-TextValue $ws "D3" "1.895.11"
Set-TextValue $ws "E3" "  +0.00%  "

# Row 4
Set-TextValue $ws "E4" "  -0.02%  "

# Row 5
Set-TextValue $ws "D5" "0.7723"
Set-TextValue $ws "E5" "  -2.33%  "

# Row 6
Set-TextValue $ws "D6" "244.52"
Set-TextValue $ws "E6" "  +0.38%  "

# Row 7
Set-TextValue $ws "E7" "  -0.01%  "

# Row 8
Set-TextValue $ws "D8" "0.3127"
Set-TextValue $ws "E8" "  -0.80%  "

# Row 9
Set-TextValue $ws "D9" "25.67"
Set-TextValue $ws "E9" "  +1.18%  "

# Row 10
Set-TextValue $ws "D10" "0.07251"

# Row 11
Set-TextValue $ws "D11" "0.08934"
Set-TextValue $ws "E11" "  +10.24%  "

# Row 12
Set-TextValue $ws "D12" "0.7720"
Set-TextValue $ws "E12" "  +0.80%  "

# Row 13
Set-TextValue $ws "D13" "5.426"
Set-TextValue $ws "E13" "  -2.65%  "

# Row 14
Set-TextValue $ws "D14" "1.877.42"
Set-TextValue $ws "E14" "  -2.12%  "

# Row 15
Set-TextValue $ws "E15" "  +2.09%  "

# Row 16
Set-TextValue $ws "D16" "6.180"
Set-TextValue $ws "E16" "  +0.14%  "

# Row 17
Set-TextValue $ws "D17" "29.863.08"
Set-TextValue $ws "E17" "  -0.17%  "

# Row 18
Set-TextValue $ws "D18" "13.94"
Set-TextValue $ws "E18" "  +0.04%  "

# Row 19
Set-TextValue $ws "D19" "245.68"
Set-TextValue $ws "E19" "  +0.71%  "

# Row 20
Set-TextValue $ws "D20" "0.000007865"
Set-TextValue $ws "E20" "  +1.08%  "

# Row 21
Set-TextValue $ws "B21" "Dai"
Set-TextValue $ws "C21" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D21" "1.000"
Set-TextValue $ws "E21" "  -0.11%  "

# Row 22
Set-TextValue $ws "B22" "Chainlink"
Set-TextValue $ws "C22" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D22" "8.125"
Set-TextValue $ws "E22" "  -1.28%  "

# Row 23
Set-TextValue $ws "D23" "2.119.69"
Set-TextValue $ws "E23" "  -0.98%  "

# Row 24
Set-TextValue $ws "E24" "  -0.09%  "

# Row 25
Set-TextValue $ws "D25" "0.1584"
Set-TextValue $ws "E25" "  -4.27%  "

# Row 26
Set-TextValue $ws "D26" "9.517"
Set-TextValue $ws "E26" "  +1.13%  "

# Row 27
Set-TextValue $ws "D27" "162.70"
Set-TextValue $ws "E27" "  -0.66%  "

# Row 28
Set-TextValue $ws "D28" "18.80"
Set-TextValue $ws "E28" "  +0.53%  "

# Row 29
Set-TextValue $ws "D29" "2.040"
Set-TextValue $ws "E29" "  -1.00%  "

# Row 30
Set-TextValue $ws "E30" "  +1.89%  "

# Row 31
Set-TextValue $ws "B31" "PancakeSwap"
Set-TextValue $ws "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D31" "1.544"
Set-TextValue $ws "E31" "  -0.35%  "

# Row 32
Set-TextValue $ws "B32" "Filecoin"
Set-TextValue $ws "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D32" "4.544"
Set-TextValue $ws "E32" "  +1.65%  "

# Row 33
Set-TextValue $ws "E33" "  +0.37%  "

# Row 34
Set-TextValue $ws "D34" "0.05500"
Set-TextValue $ws "E34" "  -0.51%  "

# Row 35
Set-TextValue $ws "E35" "  -2.14%  "

# Row 36
Set-TextValue $ws "D36" "0.7490"
Set-TextValue $ws "E36" "  +1.33%  "

# Row 37
Set-TextValue $ws "D37" "0.9990"
Set-TextValue $ws "E37" "  -0.26%  "

# Row 38
Set-TextValue $ws "D38" "2.710"
Set-TextValue $ws "E38" "  +3.44%  "

# Row 39
Set-TextValue $ws "D39" "0.01958"
Set-TextValue $ws "E39" "  +1.72%  "

# Row 40
Set-TextValue $ws "D40" "2.787"
Set-TextValue $ws "E40" "  +0.29%  "

# Row 41
Set-TextValue $ws "D41" "0.4501"
Set-TextValue $ws "E41" "  +1.69%  "

# Row 42
Set-TextValue $ws "D42" "73.75"
Set-TextValue $ws "E42" "  -0.92%  "

# Row 43
Set-TextValue $ws "B43" "Maker"
Set-TextValue $ws "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D43" "1.090.93"
Set-TextValue $ws "E43" "  -5.19%  "

# Row 44
Set-TextValue $ws "B44" "FraxShare"
Set-TextValue $ws "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D44" "6.030"
Set-TextValue $ws "E44" "  +2.48%  "

# Row 45
Set-TextValue $ws "D45" "0.8544"
Set-TextValue $ws "E45" "  +0.20%  "

# Row 46
Set-TextValue $ws "E46" "  +0.00%  "

# Row 47
Set-TextValue $ws "D47" "1.882"
Set-TextValue $ws "E47" "  +0.26%  "

# Row 48
Set-TextValue $ws "D48" "102.46"
Set-TextValue $ws "E48" "  -2.06%  "

# Row 49
Set-TextValue $ws "D49" "7.602"
Set-TextValue $ws "E49" "  +1.96%  "

# Row 50
Set-TextValue $ws "D50" "9.835"
Set-TextValue $ws "E50" "  -2.32%  "

# Row 51
Set-TextValue $ws "D51" "2.976"
Set-TextValue $ws "E51" "  -2.10%  "
